$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The scrape re-run dropped the "Naturaline Damen Bustier Weiss L" product
# (previously row 3) from the result set - delete that row, shifting every
# following product row up by one.
$ws.Rows.Item(3).Delete()

# Every remaining data row (header stays at row 1) was re-scraped at the
# new run's timestamp, so stamp column O (timestamp) on rows 2-71.
for ($r = 2; $r -le 71; $r++) {
    $ws.Cells.Item($r, 15).Value = "2022-08-28 20:57:25"
}
